$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 100, pushing existing rows 100-198 down to 101-199.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new record's data.
$ws.Cells.Item(100, 1).Value = 11
$ws.Cells.Item(100, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(100, 3).Value = "Bíobío"
$ws.Cells.Item(100, 4).Value = 44789
$ws.Cells.Item(100, 5).Value = 8
$ws.Cells.Item(100, 6).Value = 100112003
$ws.Cells.Item(100, 7).Value = "Ajo"
$ws.Cells.Item(100, 8).Value = "Chino"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 400
$ws.Cells.Item(100, 11).Value = 23000
$ws.Cells.Item(100, 12).Value = 24000
$ws.Cells.Item(100, 13).Value = 23500
$ws.Cells.Item(100, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(100, 15).Value = "China"
$ws.Cells.Item(100, 16).Value = 2350
$ws.Cells.Item(100, 17).Value = 10
$ws.Cells.Item(100, 18).Value = "Hortaliza"
